$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D6").Value = "maa://42407 (96.67)"
$ws.Range("D7").Value = "maa://21955 (95.12)"
$ws.Range("L7").Value = "maa://28624 (92.92), maa://24957 (97.73)"
$ws.Range("A8").Value = "更新日期：2025.03.20 13:18:45"
$ws.Range("AF8").Value = "*maa://24479 (78.65), *maa://21990 (51.85)"
$ws.Range("P15").Value = "maa://24762 (90.48), *maa://22727 (70.0)"
$ws.Range("D17").Value = "maa://21624 (85.0)"
$ws.Range("AF18").Value = "*maa://24313 (59.52), **maa://29784 (48.28), maa://47854 (100.0)"
$ws.Range("AB19").Value = "*maa://30709 (66.22), *maa://36668 (57.5)"
$ws.Range("AB21").Value = "maa://21443 (81.41), ***maa://23820 (30.0)"
$ws.Range("D23").Value = "***maa://28036 (28.38), *maa://41753 (52.38)"
$ws.Range("L23").Value = "maa://39756 (95.91), maa://39875 (94.59)"
$ws.Range("D24").Value = "*maa://24368 (78.44), *maa://46650 (62.5)"
$ws.Range("X24").Value = "maa://29988 (84.17), maa://23504 (93.33), **maa://22892 (40.14), *maa://25141 (77.1), *maa://36663 (77.5), ***maa://22815 (23.08)"
$ws.Range("AF28").Value = "maa://36660 (92.54), *maa://36701 (66.67)"
$ws.Range("AF29").Value = "*maa://24080 (68.85), maa://42865 (82.09), ***maa://34960 (8.33)"
$ws.Range("L31").Value = "maa://35926 (93.36), maa://36258 (84.87), *maa://43904 (72.73)"
$ws.Range("L32").Value = "maa://28065 (95.56)"
$ws.Range("L35").Value = "maa://41296 (96.45)"
$ws.Range("L37").Value = "maa://45718 (98.18), *maa://47069 (73.33), maa://45789 (100.0)"
$ws.Range("H39").Value = "maa://36670 (89.22), maa://25199 (84.82), maa://30434 (91.86), maa://45059 (80.95), ***maa://25036 (19.23), *maa://44165 (66.67)"
$ws.Range("T39").Value = "maa://45788 (80.77), maa://47079 (92.31), *maa://45790 (73.33)"
$ws.Range("P49").Value = "*maa://39643 (64.71)"
$ws.Range("H55").Value = "maa://32532 (92.04)"
$ws.Range("H59").Value = "maa://31270 (94.74), maa://27746 (82.46)"

$wb.Save()
